$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rotate CS_*_HT labels in AW1:BC1 left by one ---
$ws.Range("AW1").Value = "Odd_CS_0-1_HT"
$ws.Range("AX1").Value = "Odd_CS_0-2_HT"
$ws.Range("AY1").Value = "Odd_CS_1-2_HT"
$ws.Range("AZ1").Value = "Odd_CS_0-3_HT"
$ws.Range("BA1").Value = "Odd_CS_1-3_HT"
$ws.Range("BB1").Value = "Odd_CS_2-3_HT"
$ws.Range("BC1").Value = "Odd_CS_3-3_HT"

# --- Row 2: update match data ---
$ws.Range("A2").Value = "Yq39Z5Qj"
$ws.Range("B2").Value = "29/10/2024"
$ws.Range("C2").Value = "14:00"
$ws.Range("D2").Value = "SLOVAKIA - NIKE LIGA"
$ws.Range("E2").Value = "Dun. Streda"
$ws.Range("F2").Value = "Podbrezova"
$ws.Range("G2").Value = 1.7
$ws.Range("H2").Value = 3.85
$ws.Range("I2").Value = 4.6
$ws.Range("J2").Value = 2.15
$ws.Range("K2").Value = 2.35
$ws.Range("L2").Value = 4.65
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 8.5
$ws.Range("O2").Value = 1.26
$ws.Range("P2").Value = 3.65
$ws.Range("Q2").Value = 1.78
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 1.32
$ws.Range("T2").Value = 3.25
$ws.Range("U2").Value = 1.75
$ws.Range("V2").Value = 1.95
$ws.Range("W2").Value = 7.1
$ws.Range("X2").Value = 8.75
$ws.Range("Y2").Value = 8.5
$ws.Range("Z2").Value = 14
$ws.Range("AA2").Value = 14
$ws.Range("AB2").Value = 26
$ws.Range("AC2").Value = 8.5
$ws.Range("AD2").Value = 7.8
$ws.Range("AE2").Value = 16.5
$ws.Range("AF2").Value = 75
$ws.Range("AG2").Value = 600
$ws.Range("AH2").Value = 13
$ws.Range("AI2").Value = 30
$ws.Range("AJ2").Value = 16
$ws.Range("AK2").Value = 90
$ws.Range("AL2").Value = 45
$ws.Range("AM2").Value = 50
$ws.Range("AN2").Value = 3.7
$ws.Range("AO2").Value = 7.8
$ws.Range("AP2").Value = 15
$ws.Range("AQ2").Value = 24
$ws.Range("AR2").Value = 45
$ws.Range("AS2").Value = 175
$ws.Range("AT2").Value = 3.25
$ws.Range("AU2").Value = 7
$ws.Range("AV2").Value = 55
$ws.Range("AW2").Value = 6.5
$ws.Range("AX2").Value = 24
$ws.Range("AY2").Value = 26
$ws.Range("AZ2").Value = 120
$ws.Range("BA2").Value = 150
$ws.Range("BB2").Value = 300
$ws.Range("BC2").Value = 51
$ws.Range("BD2").Value = 51
$ws.Range("A3").Value = "SSDbN2yR"
$ws.Range("B3").Value = "29/10/2024"
$ws.Range("C3").Value = "14:00"
$ws.Range("D3").Value = "SLOVAKIA - NIKE LIGA"
$ws.Range("E3").Value = "Zilina"
$ws.Range("F3").Value = "Ruzomberok"
$ws.Range("G3").Value = 1.42
$ws.Range("H3").Value = 4.45
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 1.87
$ws.Range("K3").Value = 2.45
$ws.Range("L3").Value = 6.1
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 9.25
$ws.Range("O3").Value = 1.21
$ws.Range("P3").Value = 4.05
$ws.Range("Q3").Value = 1.65
$ws.Range("R3").Value = 2.18
$ws.Range("S3").Value = 1.31
$ws.Range("T3").Value = 3.3
$ws.Range("U3").Value = 1.83
$ws.Range("V3").Value = 1.87
$ws.Range("W3").Value = 7
$ws.Range("X3").Value = 7.4
$ws.Range("Y3").Value = 8.75
$ws.Range("Z3").Value = 10.25
$ws.Range("AA3").Value = 12
$ws.Range("AB3").Value = 27
$ws.Range("AC3").Value = 9.25
$ws.Range("AD3").Value = 9.25
$ws.Range("AE3").Value = 19
$ws.Range("AF3").Value = 90
$ws.Range("AG3").Value = 700
$ws.Range("AH3").Value = 17.5
$ws.Range("AI3").Value = 50
$ws.Range("AJ3").Value = 23
$ws.Range("AK3").Value = 175
$ws.Range("AL3").Value = 80
$ws.Range("AM3").Value = 65
$ws.Range("AN3").Value = 3.35
$ws.Range("AO3").Value = 6.4
$ws.Range("AP3").Value = 15
$ws.Range("AQ3").Value = 17.5
$ws.Range("AR3").Value = 40
$ws.Range("AS3").Value = 175
$ws.Range("AT3").Value = 3.3
$ws.Range("AU3").Value = 7.7
$ws.Range("AV3").Value = 60
$ws.Range("AW3").Value = 8.25
$ws.Range("AX3").Value = 37
$ws.Range("AY3").Value = 35
$ws.Range("AZ3").Value = 250
$ws.Range("BA3").Value = 250
$ws.Range("BB3").Value = 400
$ws.Range("BC3").Value = 51
$ws.Range("BD3").Value = 51

Write-Host "Edit applied successfully"
